$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("A8").Value = " Solving bug on starting the game."
